# Update column F (dSF) values on Sheet1 per repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 5
    3  = 5
    4  = -4
    5  = 2
    6  = 1
    7  = 5
    8  = 2
    9  = -2
    10 = 2
    11 = -2
    13 = 4
    14 = -1
    15 = -1
    16 = -2
    17 = 1
    18 = -2
    19 = -2
    20 = 7
    21 = 1
    22 = 1
    23 = 2
    24 = -1
    25 = -1
    26 = 1
    27 = -1
    28 = -2
    30 = 1
    31 = 2
    33 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
